$d = $word.ActiveDocument

# --- Target cell: Table 1, Row 16, Column 3 (the empty "-0.25" score cell
# for the "Dang nhap voi Facebook" row) ---
$t = $d.Tables.Item(1)
$cell = $t.Cell(16, 3)
$p = $cell.Range.Paragraphs.Item(1)

# Collapse to the zero-width point right before the cell's end-of-paragraph
# mark, and plant the _GoBack bookmark there first. Adding a bookmark with
# an existing name moves it (the old _GoBack after "Phieu cham diem" is
# removed automatically).
$bmRange = $p.Range.Duplicate
[void]$bmRange.MoveEnd(1, -1)
$bmRange.Collapse(0)
$bmRange.Bookmarks.Add("_GoBack")

# Insert the new run's text before that same point; since the bookmark
# degenerates to a point there, the inserted text lands before it, giving
# <w:r>-0.25</w:r><w:bookmarkStart/><w:bookmarkEnd/>.
$cell.Range.InsertBefore("-0.25")
